$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# B7 previously held "Comprobar que funciona y eliminar equipo si no tiene jugadores".
# Change it to reuse the existing text "Comprobar que funciona" (same as B15/B17).
$ws.Range("B7").Value = "Comprobar que funciona"

# C10 was empty; it should now contain the (new) text that replaces the old B7 string slot.
$ws.Range("C10").Value = "Eliminar organization si no tiene equipo"

# Update the active selection to B14.
$ws.Range("B14").Select()

$wb.Save()
